# Refresh the cryptos snapshot: update Price (D) / Volume(1h) (E) columns for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.549.96"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").Value = "1.916.48"
$ws.Range("E3").Value = "  +2.81%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'315.58"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "'0.5161"
$ws.Range("E7").Value = "  +3.76%  "
$ws.Range("D8").Value = "'0.3990"
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").Value = "'0.09853"
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("E10").Value = "  +2.86%  "
$ws.Range("D11").Value = "'42.31"
$ws.Range("E11").Value = "  +2.49%  "
$ws.Range("D12").Value = "'6.523"
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("E13").Value = "  +1.50%  "
$ws.Range("D14").Value = "1.922.53"
$ws.Range("E14").Value = "  +3.42%  "
$ws.Range("D15").Value = "'7.474"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").Value = "'0.00001138"
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").Value = "'94.64"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("D19").Value = "'0.06662"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("E20").Value = "  +5.09%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "'6.314"
$ws.Range("E22").Value = "  +4.04%  "
$ws.Range("D23").Value = "28.607.54"
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("D25").Value = "'2.320"
$ws.Range("E25").Value = "  +3.12%  "
$ws.Range("D26").Value = "'2.683"
$ws.Range("E26").Value = "  +8.42%  "
$ws.Range("D27").Value = "2.134.85"
$ws.Range("E27").Value = "  +2.93%  "
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("D29").Value = "'157.66"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").Value = "'129.66"
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("D31").Value = "'1.115"
$ws.Range("E31").Value = "  +6.17%  "
$ws.Range("E32").Value = "  +1.79%  "
$ws.Range("D33").Value = "'5.730"
$ws.Range("E33").Value = "  +1.48%  "
$ws.Range("D34").Value = "'3.632"
$ws.Range("E34").Value = "  +1.22%  "
$ws.Range("D35").Value = "'9.840"
$ws.Range("E35").Value = "  +6.24%  "
$ws.Range("D36").Value = "'0.06768"
$ws.Range("E36").Value = "  -0.50%  "
$ws.Range("D37").Value = "'0.02441"
$ws.Range("E37").Value = "  +2.47%  "
$ws.Range("E38").Value = "  +5.37%  "
$ws.Range("D39").Value = "'0.2229"
$ws.Range("E39").Value = "  +2.90%  "
$ws.Range("D40").Value = "'11.82"
$ws.Range("E40").Value = "  +2.82%  "
$ws.Range("D41").Value = "'0.6490"
$ws.Range("D42").Value = "'5.093"
$ws.Range("E42").Value = "  +1.45%  "
$ws.Range("E43").Value = "  +0.91%  "
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "'13.61"
$ws.Range("E45").Value = "  +1.97%  "
$ws.Range("D46").Value = "'0.6107"
$ws.Range("E46").Value = "  +2.01%  "
$ws.Range("D47").Value = "'3.775"
$ws.Range("E47").Value = "  +2.39%  "
$ws.Range("D48").Value = "'1.288"
$ws.Range("E48").Value = "  +0.51%  "
$ws.Range("D49").Value = "'2.071"
$ws.Range("E49").Value = "  +4.68%  "
$ws.Range("D50").Value = "'124.87"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("E51").Value = "  +1.11%  "
